$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 875
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H103").Value = 1389.7273
$ws.Range("I103").Value = 844
$ws.Range("J103").Value = 1511
$ws.Range("K103").Value = 2532
$ws.Range("L103").Value = 4533
$ws.Range("M103").Value = -1946
$ws.Range("N103").Value = -5705

$ws.Range("H106").Value = 3166.3333
$ws.Range("I106").Value = 3166.3333
$ws.Range("K106").Value = 3166.3333
$ws.Range("M106").Value = -2535.3333

$ws.Range("H138").Value = 3341.0557
$ws.Range("J138").Value = 4995.4443
$ws.Range("L138").Value = 14986.3329
$ws.Range("N138").Value = -25266.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4165.3335
$ws.Range("I132").Value = 3748.5
$ws.Range("K132").Value = 11245.5
$ws.Range("M132").Value = -8715.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 100000
$ws.Range("J9").Value = 100000
$ws.Range("L9").Value = 100000
$ws.Range("N9").Value = -100336

$ws.Range("H44").Value = 60000
$ws.Range("J44").Value = 60000
$ws.Range("L44").Value = 60000
$ws.Range("N44").Value = -60994

$ws.Range("H94").Value = 3832.75
$ws.Range("I94").Value = 2237.4285
$ws.Range("K94").Value = 2237.4285
$ws.Range("M94").Value = -1786.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376

$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880

$ws.Range("H99").Value = 3465.25
$ws.Range("I99").Value = 3465.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3465.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1967.25
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 3465.25
$ws.Range("I126").Value = 3465.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10395.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7925.75
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 14187.5
$ws.Range("I80").Value = 9000
$ws.Range("J80").Value = 14928.571
$ws.Range("K80").Value = 27000
$ws.Range("L80").Value = 44785.713
$ws.Range("M80").Value = -26064
$ws.Range("N80").Value = -46657.713

$ws.Range("H83").Value = 14187.5
$ws.Range("I83").Value = 9000
$ws.Range("J83").Value = 14928.571
$ws.Range("K83").Value = 81000
$ws.Range("L83").Value = 134357.139
$ws.Range("M83").Value = -76320
$ws.Range("N83").Value = -143717.139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1298.3334
$ws.Range("I16").Value = 1298.3334
$ws.Range("K16").Value = 1298.3334
$ws.Range("M16").Value = -1128.3334

$ws.Range("H22").Value = 7472.9165
$ws.Range("I22").Value = 9655.556
$ws.Range("J22").Value = 925
$ws.Range("K22").Value = 9655.556
$ws.Range("L22").Value = 925
$ws.Range("M22").Value = -9360.556
$ws.Range("N22").Value = -1515

$ws.Range("H27").Value = 7472.9165
$ws.Range("I27").Value = 9655.556
$ws.Range("J27").Value = 925
$ws.Range("K27").Value = 9655.556
$ws.Range("L27").Value = 925
$ws.Range("M27").Value = -9548.556
$ws.Range("N27").Value = -1139

$ws.Range("H55").Value = 5500
$ws.Range("I55").Value = 5500
$ws.Range("K55").Value = 5500
$ws.Range("M55").Value = -5327

$ws.Range("H68").Value = 1799
$ws.Range("I68").Value = 2263.5
$ws.Range("J68").Value = 870
$ws.Range("K68").Value = 2263.5
$ws.Range("L68").Value = 870
$ws.Range("M68").Value = -1514.5
$ws.Range("N68").Value = -2368

$ws.Range("H71").Value = 1799
$ws.Range("I71").Value = 2263.5
$ws.Range("J71").Value = 870
$ws.Range("K71").Value = 11317.5
$ws.Range("L71").Value = 4350
$ws.Range("M71").Value = -7573.5
$ws.Range("N71").Value = -11838

$ws.Range("H82").Value = 3232.2856
$ws.Range("I82").Value = 2971
$ws.Range("J82").Value = 4800
$ws.Range("K82").Value = 2971
$ws.Range("L82").Value = 4800
$ws.Range("M82").Value = -2610
$ws.Range("N82").Value = -5522

$ws.Range("H85").Value = 3232.2856
$ws.Range("I85").Value = 2971
$ws.Range("J85").Value = 4800
$ws.Range("K85").Value = 2971
$ws.Range("L85").Value = 4800
$ws.Range("M85").Value = -1723
$ws.Range("N85").Value = -7296

$ws.Range("H93").Value = 1036.25
$ws.Range("J93").Value = 1066.3334
$ws.Range("L93").Value = 1066.3334
$ws.Range("N93").Value = -3562.3334

$ws.Range("H122").Value = 5197.4
$ws.Range("I122").Value = 3997
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 11991
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -9541
$ws.Range("N122").Value = -34897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2500250
$ws.Range("J3").Value = 500
$ws.Range("L3").Value = 500
$ws.Range("N3").Value = -728

$ws.Range("H62").Value = 5627.4546
$ws.Range("I62").Value = 2438
$ws.Range("J62").Value = 7450
$ws.Range("K62").Value = 2438
$ws.Range("L62").Value = 7450
$ws.Range("M62").Value = -1814
$ws.Range("N62").Value = -8698

$ws.Range("H65").Value = 5627.4546
$ws.Range("I65").Value = 2438
$ws.Range("J65").Value = 7450
$ws.Range("K65").Value = 12190
$ws.Range("L65").Value = 37250
$ws.Range("M65").Value = -9070
$ws.Range("N65").Value = -43490

$ws.Range("H107").Value = 899.1111
$ws.Range("I107").Value = 865.1667
$ws.Range("J107").Value = 967
$ws.Range("K107").Value = 2595.5001
$ws.Range("L107").Value = 2901
$ws.Range("M107").Value = -675.5001000000002
$ws.Range("N107").Value = -6741

$ws.Range("H136").Value = 3144.1904
$ws.Range("I136").Value = 3212
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 9636
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -7086
$ws.Range("N136").Value = -12600
